$d = $word.ActiveDocument

$d.Content.Find.Execute("803×9=", $true, $false, $false, $false, $false, $true, 1, $false, "367×5=", 2)
$d.Content.Find.Execute("953×7=", $true, $false, $false, $false, $false, $true, 1, $false, "215×2=", 2)
$d.Content.Find.Execute("270×4=", $true, $false, $false, $false, $false, $true, 1, $false, "520×9=", 2)
$d.Content.Find.Execute("823×3=", $true, $false, $false, $false, $false, $true, 1, $false, "816×2=", 2)
$d.Content.Find.Execute("349×8=", $true, $false, $false, $false, $false, $true, 1, $false, "484×3=", 2)
$d.Content.Find.Execute("486×6=", $true, $false, $false, $false, $false, $true, 1, $false, "850×5=", 2)
$d.Content.Find.Execute("932×5=", $true, $false, $false, $false, $false, $true, 1, $false, "346×2=", 2)
$d.Content.Find.Execute("104×7=", $true, $false, $false, $false, $false, $true, 1, $false, "665×5=", 2)
$d.Content.Find.Execute("584×7=", $true, $false, $false, $false, $false, $true, 1, $false, "484×9=", 2)
$d.Content.Find.Execute("453×5=", $true, $false, $false, $false, $false, $true, 1, $false, "969×4=", 2)
$d.Content.Find.Execute("270×8=", $true, $false, $false, $false, $false, $true, 1, $false, "584×3=", 2)
$d.Content.Find.Execute("204×8=", $true, $false, $false, $false, $false, $true, 1, $false, "107×4=", 2)
$d.Content.Find.Execute("747×5=", $true, $false, $false, $false, $false, $true, 1, $false, "362×9=", 2)
$d.Content.Find.Execute("166×6=", $true, $false, $false, $false, $false, $true, 1, $false, "154×3=", 2)
$d.Content.Find.Execute("293×8=", $true, $false, $false, $false, $false, $true, 1, $false, "130×2=", 2)
$d.Content.Find.Execute("974×5=", $true, $false, $false, $false, $false, $true, 1, $false, "788×7=", 2)
$d.Content.Find.Execute("734×7=", $true, $false, $false, $false, $false, $true, 1, $false, "869×9=", 2)
$d.Content.Find.Execute("926×9=", $true, $false, $false, $false, $false, $true, 1, $false, "629×4=", 2)
$d.Content.Find.Execute("526×2=", $true, $false, $false, $false, $false, $true, 1, $false, "605×7=", 2)
$d.Content.Find.Execute("493×9=", $true, $false, $false, $false, $false, $true, 1, $false, "834×8=", 2)
$d.Content.Find.Execute("531×7=", $true, $false, $false, $false, $false, $true, 1, $false, "834×3=", 2)
$d.Content.Find.Execute("951×7=", $true, $false, $false, $false, $false, $true, 1, $false, "142×3=", 2)
$d.Content.Find.Execute("568×9=", $true, $false, $false, $false, $false, $true, 1, $false, "849×6=", 2)
$d.Content.Find.Execute("934×5=", $true, $false, $false, $false, $false, $true, 1, $false, "235×7=", 2)
$d.Content.Find.Execute("131×5=", $true, $false, $false, $false, $false, $true, 1, $false, "292×9=", 2)
